$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 149
$ws.Range("B149").Value = 63902
$ws.Range("E149").Value = 34.04
$ws.Range("F149").Value = 2
$ws.Range("G149").Value = 64.04000000000001
# Row 150
$ws.Range("B150").Value = 48654
$ws.Range("E150").Value = 38.26
$ws.Range("F150").Value = -1
$ws.Range("G150").Value = -32.02
# Row 257
$ws.Range("B257").Value = 65289
$ws.Range("C257").Value = 'HIM-ALMOND &amp;amp; ROSE SOAP 125G'
$ws.Range("E257").Value = 34.09
# Row 258
$ws.Range("B258").Value = 65317
$ws.Range("E258").Value = 213.75
# Row 263
$ws.Range("B263").Value = 65304
$ws.Range("E263").Value = 511.63
# Row 264
$ws.Range("B264").Value = 64979
$ws.Range("E264").Value = 314.41
$ws.Range("F264").Value = 82
$ws.Range("G264").Value = 24251.5
# Row 265
$ws.Range("B265").Value = 48719
$ws.Range("E265").Value = 353.35
$ws.Range("F265").Value = -81
$ws.Range("G265").Value = -23955.75
# Row 266
$ws.Range("B266").Value = 65297
$ws.Range("E266").Value = 96.34999999999999
# Row 269
$ws.Range("B269").Value = 65299
$ws.Range("E269").Value = 71.66
# Row 270
$ws.Range("B270").Value = 65300
$ws.Range("E270").Value = 125.37
# Row 271
$ws.Range("B271").Value = 65308
$ws.Range("E271").Value = 97.98
# Row 274
$ws.Range("B274").Value = 65290
$ws.Range("C274").Value = 'HIM-CUCUMBER &amp;amp; COCONUT SOAP 125G'
$ws.Range("E274").Value = 34.09
# Row 276
$ws.Range("B276").Value = 65282
$ws.Range("E276").Value = 35.4
# Row 281
$ws.Range("B281").Value = 65301
$ws.Range("E281").Value = 102.41
# Row 282
$ws.Range("B282").Value = 65305
$ws.Range("E282").Value = 58.16
# Row 286
$ws.Range("B286").Value = 65320
$ws.Range("E286").Value = 289.3
# Row 289
$ws.Range("B289").Value = 65291
$ws.Range("C289").Value = 'HIM-HONEY &amp;amp; CREAM SOAP 125G'
$ws.Range("E289").Value = 34.09
# Row 290
$ws.Range("B290").Value = 65288
$ws.Range("C290").Value = 'HIM-HONEY &amp;amp; CREAM SOAP 75G'
$ws.Range("E290").Value = 21.5
# Row 292
$ws.Range("B292").Value = 65319
$ws.Range("E292").Value = 77.73
# Row 296
$ws.Range("B296").Value = 65307
$ws.Range("E296").Value = 400.95
# Row 297
$ws.Range("B297").Value = 65280
$ws.Range("E297").Value = 49.82
# Row 299
$ws.Range("B299").Value = 65314
$ws.Range("E299").Value = 32.13
# Row 303
$ws.Range("B303").Value = 65313
$ws.Range("E303").Value = 92.90000000000001
# Row 305
$ws.Range("B305").Value = 65306
$ws.Range("E305").Value = 107.59
# Row 306
$ws.Range("B306").Value = 65303
$ws.Range("E306").Value = 496.69
# Row 309
$ws.Range("B309").Value = 65302
$ws.Range("E309").Value = 449.92
# Row 313
$ws.Range("B313").Value = 62997
$ws.Range("F313").Value = 72
$ws.Range("G313").Value = 22020.48
# Row 314
$ws.Range("B314").Value = 57854
$ws.Range("F314").Value = 2
$ws.Range("G314").Value = 611.6799999999999
# Row 316
$ws.Range("B316").Value = 63565
$ws.Range("D316").Value = 102.71
$ws.Range("E316").Value = 109.19
$ws.Range("F316").Value = 60
$ws.Range("G316").Value = 6162.6
# Row 317
$ws.Range("B317").Value = 57077
$ws.Range("D317").Value = 93.08
$ws.Range("E317").Value = 111.2
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 93.08
# Row 318
$ws.Range("B318").Value = 61610
$ws.Range("E318").Value = 122.71
$ws.Range("F318").Value = -58
$ws.Range("G318").Value = -5957.18
# Row 351
$ws.Range("B351").Value = 63531
$ws.Range("F351").Value = 80
$ws.Range("G351").Value = 11478.4
# Row 352
$ws.Range("B352").Value = 63571
$ws.Range("F352").Value = 29
$ws.Range("G352").Value = 4160.92
# Row 382
$ws.Range("B382").Value = 60325
$ws.Range("E382").Value = 151.57
$ws.Range("F382").Value = -102
$ws.Range("G382").Value = -12939.72
# Row 383
$ws.Range("B383").Value = 63560
$ws.Range("E383").Value = 134.87
$ws.Range("F383").Value = 104
$ws.Range("G383").Value = 13193.44
# Row 421
$ws.Range("B421").Value = 63008
$ws.Range("F421").Value = 504
$ws.Range("G421").Value = 76189.67999999999
# Row 422
$ws.Range("B422").Value = 57857
$ws.Range("F422").Value = 3
$ws.Range("G422").Value = 453.51
# Row 431
$ws.Range("B431").Value = 63102
$ws.Range("C431").Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Range("F431").Value = 36
$ws.Range("G431").Value = 2140.92
# Row 432
$ws.Range("B432").Value = 53082
$ws.Range("C432").Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Range("F432").Value = 1
$ws.Range("G432").Value = 59.47
# Row 536
$ws.Range("B536").Value = 58047
$ws.Range("D536").Value = 105.54
$ws.Range("E536").Value = 126.1
$ws.Range("F536").Value = 54
$ws.Range("G536").Value = 5699.16
# Row 537
$ws.Range("B537").Value = 47097
$ws.Range("D537").Value = 112.28
$ws.Range("E537").Value = 134.16
$ws.Range("F537").Value = 15
$ws.Range("G537").Value = 1684.2
# Row 579
$ws.Range("B579").Value = 53757
$ws.Range("E579").Value = 16.08
$ws.Range("F579").Value = -159
$ws.Range("G579").Value = -2138.55
# Row 580
$ws.Range("B580").Value = 65069
$ws.Range("E580").Value = 14.3
$ws.Range("F580").Value = 172
$ws.Range("G580").Value = 2313.4
# Row 583
$ws.Range("B583").Value = 65066
$ws.Range("E583").Value = 13.61
$ws.Range("F583").Value = 313
$ws.Range("G583").Value = 4009.53
# Row 584
$ws.Range("B584").Value = 53263
$ws.Range("E584").Value = 15.29
$ws.Range("F584").Value = -309
$ws.Range("G584").Value = -3958.29
# Row 586
$ws.Range("B586").Value = 45695
$ws.Range("E586").Value = 23.58
$ws.Range("F586").Value = -36
$ws.Range("G586").Value = -710.28
# Row 587
$ws.Range("B587").Value = 64915
$ws.Range("E587").Value = 20.98
$ws.Range("F587").Value = 40
$ws.Range("G587").Value = 789.2
# Row 599
$ws.Range("B599").Value = 45709
$ws.Range("E599").Value = 15.69
$ws.Range("F599").Value = -300
$ws.Range("G599").Value = -3945
# Row 600
$ws.Range("B600").Value = 64925
$ws.Range("E600").Value = 13.97
$ws.Range("F600").Value = 302
$ws.Range("G600").Value = 3971.3
# Row 601
$ws.Range("B601").Value = 45702
$ws.Range("E601").Value = 31.43
$ws.Range("F601").Value = -215
$ws.Range("G601").Value = -5654.5
# Row 602
$ws.Range("B602").Value = 64919
$ws.Range("E602").Value = 27.97
$ws.Range("F602").Value = 224
$ws.Range("G602").Value = 5891.2
# Row 687
$ws.Range("B687").Value = 53319
$ws.Range("E687").Value = 310.64
$ws.Range("F687").Value = -6
$ws.Range("G687").Value = -1643.52
# Row 688
$ws.Range("B688").Value = 64810
$ws.Range("E688").Value = 291.22
$ws.Range("F688").Value = 7
$ws.Range("G688").Value = 1917.44
# Row 720
$ws.Range("B720").Value = 60022
$ws.Range("E720").Value = 37.22
$ws.Range("F720").Value = -113
$ws.Range("G720").Value = -3709.79
# Row 721
$ws.Range("B721").Value = 64830
$ws.Range("E721").Value = 34.9
$ws.Range("F721").Value = 117
$ws.Range("G721").Value = 3841.11
# Row 872
$ws.Range("B872").Value = 54751
$ws.Range("E872").Value = 46.34
$ws.Range("F872").Value = -19
$ws.Range("G872").Value = -776.53
# Row 873
$ws.Range("B873").Value = 65079
$ws.Range("E873").Value = 43.44
$ws.Range("F873").Value = 21
$ws.Range("G873").Value = 858.27

Write-Host "Applied all cell updates"
